# Commit: "Automatic update of files."
#
# The sheet tracks a "Förändrad" (last-changed) date in column C for every
# record row. This automatic refresh bumps that date forward by one day
# (2026-02-22 -> 2026-02-23, i.e. Excel serial 46075 -> 46076) for every
# data row currently stamped with the old date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -eq 46075) {
        $cell.Value = 46076
    }
}
